{"js": "// Change 1: \"3. Risks to meeting fishery management objectives\" intro paragraph.\n// Swap the U.S. Caribbean region bounding-box sentence for the Gulf of\n// America placeholder version (region name + \"xx degrees\" placeholders).\nconst risksHits = context.document.body.search(\n  \"Unless otherwise specified, physical indicators reported for the U.S. Caribbean region were calculated over a bounding box with limits of longitude 68 degrees W to 64.5 degrees W and latitude 17.5 degrees N to 18.75 degrees N.\",\n  { matchCase: true }\n);\nrisksHits.load(\"text\");\nawait context.sync();\n\nif (risksHits.items.length > 0) {\n  risksHits.items[0].insertText(\n    \"Unless otherwise specified, physical indicators reported for the U.S. Gulf of America region were calculated over a bounding box with limits of longitude xx degrees W to xx degrees W and latitude xx degrees N to xx degrees N.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// Change 2: PCA / traffic-light-plot paragraph. The bold placeholder run\n// \"?@fig-traffic\" (wrapped in parentheses) becomes plain, non-bold \"(figure)\"\n// text. Rewrite the whole paragraph's text in one shot so the three runs\n// (plain / bold \"?@fig-traffic\" / plain) collapse back into a single plain\n// run, matching how Word merges identically-formatted adjacent runs.\nconst pcaParas = context.document.body.paragraphs;\npcaParas.load(\"text\");\nawait context.sync();\n\nconst pcaPara = pcaParas.items.find((p) => p.text.includes(\"?@fig-traffic\"));\nif (pcaPara) {\n  const updated = pcaPara.text.replace(\"(?@fig-traffic)\", \"(figure)\");\n  pcaPara.getRange().insertText(updated, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Change 1: \"3. Risks to meeting fishery management objectives\" intro\n# paragraph. Swap the U.S. Caribbean region bounding-box sentence for the\n# Gulf of America placeholder version (region name + \"xx degrees\" values).\n$d = $word.ActiveDocument\n\n$range1 = $d.Content\n$range1.Find.ClearFormatting()\n$range1.Find.Execute(\n    \"U.S. Caribbean region\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"U.S. Gulf of America region\",\n    2\n)\n\n$range2 = $d.Content\n$range2.Find.ClearFormatting()\n$range2.Find.Execute(\n    \"68 degrees W to 64.5 degrees W\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"xx degrees W to xx degrees W\",\n    2\n)\n\n$range3 = $d.Content\n$range3.Find.ClearFormatting()\n$range3.Find.Execute(\n    \"17.5 degrees N to 18.75 degrees N\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"xx degrees N to xx degrees N\",\n    2\n)\n\n# Change 2: PCA / traffic-light-plot paragraph. The bold placeholder run\n# \"?@fig-traffic\" (wrapped in parentheses) becomes plain, non-bold \"(figure)\"\n# text, collapsing the paragraph back down into a single run.\n$range4 = $d.Content\n$range4.Find.ClearFormatting()\n$range4.Find.Execute(\n    \"(?@fig-traffic)\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"(figure)\",\n    2\n)\n"}
